# daily auto push: 2026-01-10 04:03 UTC
# Insert a new data row for 2026/01/10 12:00 slot (row 619) into Sheet1,
# shifting all subsequent rows down by one (old row 619 -> 620, ... old row 660 -> 661).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 619 (shifts existing rows 619-660 down to 620-661).
$ws.Rows.Item(619).Insert()

# Fill in the new row's values.
# Column A holds a date formatted as plain text (e.g. "2026/01/10"); setting NumberFormat
# to Text before assigning the value (and clearing it again afterwards) keeps Excel from
# auto-converting the literal string into a date serial number, while leaving the cell's
# style untouched (same as every other cell in the column).
$ws.Range("A619").NumberFormat = "@"
$ws.Range("A619").Value = "2026/01/10"
$ws.Range("A619").ClearFormats()

$ws.Range("B619").Value = "土"
$ws.Range("C619").Value = 12
$ws.Range("D619").Value = 201
